$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the two obsolete rows (Abducted/"Похищенные"/"Уурдалган" and
#    Lost/"Утерянные"/"Жоготулган"). Both were rows 6 and 7; deleting row 6
#    twice removes both and shifts everything below up.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# 2. Update the 2019 figures for the two remaining indicators.
# ---------------------------------------------------------------------------
$ws.Range("G4").Value = 146
$ws.Range("G5").Value = 127

# ---------------------------------------------------------------------------
# 3. Add the new 2020 column (H) with its header and data values.
# ---------------------------------------------------------------------------
$ws.Range("H3").Value = 2020
$ws.Range("H4").Value = 158
$ws.Range("H5").Value = 397

# Copy the existing year-header formatting (column G) onto the new column H
# header cell so it matches the other year columns.
$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)   # xlPasteFormats

# Copy the plain data-cell formatting from the neighbouring column onto the
# new column's data cells.
$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 4. Row 5 ("Seized firearms") is now the last row of the table, so it needs
#    the heavier bottom border previously used to close off the table.
#    First extend row 2's separator formatting into the new column H, then
#    copy the whole separator row's look onto row 5.
# ---------------------------------------------------------------------------
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A2:H2").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)   # xlPasteFormats

$ws.Application.CutCopyMode = $false

$ws.Range("A1").Select()
